# New fitting and year date change
$d = $word.ActiveDocument

# --- 1. Locate the table that holds the "Szekek dobozok" rows and append
#       a new row for "KMaxxi hatso doboz jobb es baloldalra" / "55 000.- / db"
$targetTable = $null
foreach ($tbl in $d.Tables) {
    if ($tbl.Range.Text -like "*XXL hátsó doboz jobb és baldoldalra*") {
        $targetTable = $tbl
        break
    }
}

$newRow = $targetTable.Rows.Add()
$newRow.Cells(1).Range.Text = "KMaxxi hátsó doboz jobb és baloldalra"
$newRow.Cells(2).Range.Text = "55 000.- / db"
$newRow.Cells(2).Range.ParagraphFormat.Alignment = 2

# --- 2. Move the _GoBack bookmark: remove it from the old "Ráfutófék" spot
#        (Word auto-manages _GoBack at the most-recent edit, which is now
#        in the new row's price cell)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $newRow.Cells(2).Range)
